$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("E2").Value = "Maandelijkse Facturatie a €250,99I = geschatte éénmalige investeringskosten n.v.t.Y = aantal jaren looptijd overeenkomstX = geschatte jaarlijkse kosten"
$ws.Range("G2").Value = "Eelco Aartsen"
$ws.Range("I2").Value = "075-6163455"
$ws.Range("T2").Value = "Pietje Puk"
$ws.Range("V2").Value = "06-1231232"
$ws.Range("Z2").Value = "Erwtensoep 2021"

# --- Row 3 updates ---
$ws.Range("E3").Value = "Maandelijkse Facturatie a €123,34I = geschatte éénmalige investeringskosten n.v.t.Y = aantal jaren looptijd overeenkomstX = geschatte jaarlijkse kosten"
$ws.Range("G3").Value = "Eelco Aartsen"
$ws.Range("I3").Value = "075-6163455"
$ws.Range("T3").Value = "Pietje Puk"
$ws.Range("V3").Value = "06-1231232"
$ws.Range("Z3").Value = "Erwtensoep 2021"

# --- New contact-person hyperlinks (e-mail addresses) ---
$ws.Hyperlinks.Add($ws.Range("H2"), "mailto:eelco@aesset.nl", "", "", "eelco@aesset.nl")
$ws.Range("H2").Font.Color = 16711680
$ws.Range("H2").Font.Underline = $false

$ws.Hyperlinks.Add($ws.Range("U2"), "mailto:p.puk@npo.nl", "", "", "p.puk@npo.nl")
$ws.Range("U2").Font.Color = 16711680
$ws.Range("U2").Font.Underline = $false

$ws.Hyperlinks.Add($ws.Range("H3"), "mailto:eelco@aesset.nl", "", "", "eelco@aesset.nl")
$ws.Range("H3").Font.Color = 16711680
$ws.Range("H3").Font.Underline = $false

$ws.Hyperlinks.Add($ws.Range("U3"), "mailto:p.puk@npo.nl", "", "", "p.puk@npo.nl")
$ws.Range("U3").Font.Color = 16711680
$ws.Range("U3").Font.Underline = $false

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 14.15
$ws.Rows.Item(3).RowHeight = 14.15

# --- Selection / view state ---
[void]$ws.Range("T3:AA3").Select()

Write-Output "done"
